# Generate Report for Handoff
# The dc34a600-f094-46b3-8fdd-b740f256be16.md item has been re-handed-off:
#  - Status flips from "Handed back: in sync with en-US" to "Ready for handoff"
#  - The relevant timestamp columns get refreshed
#  - An error detail explaining the stale handback file is recorded (zh-cn / de-de tabs)

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b58988280f3fccd8ab251466dba89658a62e9b28/e2e/dc34a600-f094-46b3-8fdd-b740f256be16.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e0e93bc8ba1ea7d47cc57e852b5469e8ba7048f4/e2e/dc34a600-f094-46b3-8fdd-b740f256be16.md."

# --- Overview sheet: row 3 is the dc34a600-f094-46b3-8fdd-b740f256be16.md entry ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $readyForHandoff
$wsOverview.Range("F3").Value = $readyForHandoff
$wsOverview.Range("G3").Value = "2016-08-22 14:54:51"

# --- zh-cn sheet: row 3 is the dc34a600-f094-46b3-8fdd-b740f256be16.md entry ---
# ColumnWidth 39.17 (COM "characters" units) round-trips to the stored OOXML
# column width of 40 (same value already used by the other width=40 columns
# in this sheet, e.g. column A).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsZhCn.Range("C3").Value = $readyForHandoff
$wsZhCn.Range("H3").Value = "2016-08-22 14:54:46"
$wsZhCn.Range("P3").Value = $errorDetail

# --- de-de sheet: row 3 is the dc34a600-f094-46b3-8fdd-b740f256be16.md entry ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Range("C3").Value = $readyForHandoff
$wsDeDe.Range("H3").Value = "2016-08-22 14:54:51"
$wsDeDe.Range("P3").Value = $errorDetail
